# Insert a new data row at row 200 (shifting existing rows 200-237 down to
# 201-238), then populate the newly inserted row with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 200, pushing everything below
# down by one (this mirrors Excel's "Insert Sheet Rows" behaviour).
$ws.Rows("200").Insert()

# Populate the newly inserted row 200 with the new record's values.
$ws.Range("A200").Value = 11
$ws.Range("B200").Value = "Vega Monumental Concepción"
$ws.Range("C200").Value = "Bíobío"
$ws.Range("D200").Value = 44637
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = "Fruta"
$ws.Range("G200").Value = 100102
$ws.Range("H200").Value = "Cítricos"
$ws.Range("I200").Value = 100102005
$ws.Range("J200").Value = "Naranja"
$ws.Range("K200").Value = "Valencia"
$ws.Range("L200").Value = "Primera"
$ws.Range("M200").Value = 220
$ws.Range("N200").Value = 8500
$ws.Range("O200").Value = 9000
$ws.Range("P200").Value = 8727
$ws.Range("Q200").Value = "`$/caja 15 kilos granel"
$ws.Range("R200").Value = "Región de O'Higgins"
$ws.Range("S200").Value = 582
$ws.Range("T200").Value = 15
